$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2307.5
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 2076.6667
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 2076.6667
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -2426.6667
$ws.Range("H62").Value = 4665.6665
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 4665.6665
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H86").Value = 695.6667
$ws.Range("I86").Value = 655.3333
$ws.Range("K86").Value = 655.3333
$ws.Range("M86").Value = 467.6667
$ws.Range("H89").Value = 695.6667
$ws.Range("I89").Value = 655.3333
$ws.Range("K89").Value = 3276.6665
$ws.Range("M89").Value = 2339.3335
$ws.Range("H127").Value = 1198.25
$ws.Range("I127").Value = 1131.3334
$ws.Range("K127").Value = 3394.0002
$ws.Range("M127").Value = 1565.9998
$ws.Range("H132").Value = 1002.1429
$ws.Range("I132").Value = 1002.1429
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3006.4287
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -476.4287000000004
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3140.2576
$ws.Range("I32").Value = 2104.6038
$ws.Range("K32").Value = 2104.6038
$ws.Range("M32").Value = -1817.6038
$ws.Range("H45").Value = 3216835.2
$ws.Range("I45").Value = 10002410
$ws.Range("J45").Value = 2615.5789
$ws.Range("K45").Value = 10002410
$ws.Range("L45").Value = 2615.5789
$ws.Range("M45").Value = -10002033
$ws.Range("N45").Value = -3369.5789
$ws.Range("H74").Value = 1069.6072
$ws.Range("I74").Value = 487.38095
$ws.Range("J74").Value = 2816.2856
$ws.Range("K74").Value = 487.38095
$ws.Range("L74").Value = 2816.2856
$ws.Range("M74").Value = 386.61905
$ws.Range("N74").Value = -4564.2856
$ws.Range("H77").Value = 1069.6072
$ws.Range("I77").Value = 487.38095
$ws.Range("J77").Value = 2816.2856
$ws.Range("K77").Value = 2436.90475
$ws.Range("L77").Value = 14081.428
$ws.Range("M77").Value = 1931.09525
$ws.Range("N77").Value = -22817.428
$ws.Range("H132").Value = 1969.875
$ws.Range("I132").Value = 1451.5853
$ws.Range("K132").Value = 4354.7559
$ws.Range("M132").Value = -1824.7559

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31694
$ws.Range("H86").Value = 288035.44
$ws.Range("I86").Value = 2749.6
$ws.Range("J86").Value = 1001250
$ws.Range("K86").Value = 2749.6
$ws.Range("L86").Value = 1001250
$ws.Range("M86").Value = -1626.6
$ws.Range("N86").Value = -1003496
$ws.Range("H89").Value = 288035.44
$ws.Range("I89").Value = 2749.6
$ws.Range("J89").Value = 1001250
$ws.Range("K89").Value = 13748
$ws.Range("L89").Value = 5006250
$ws.Range("M89").Value = -8132
$ws.Range("N89").Value = -5017482
$ws.Range("H134").Value = 5867.32
$ws.Range("I134").Value = 6167.409
$ws.Range("K134").Value = 18502.227
$ws.Range("M134").Value = -15967.227

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 8952
$ws.Range("I5").Value = 10690
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 10690
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -10578
$ws.Range("N5").Value = -2224
$ws.Range("H7").Value = 180
$ws.Range("I7").Value = 180
$ws.Range("K7").Value = 180
$ws.Range("M7").Value = -67
$ws.Range("H31").Value = 1395.2963
$ws.Range("I31").Value = 890.75
$ws.Range("J31").Value = 1607.7368
$ws.Range("K31").Value = 890.75
$ws.Range("L31").Value = 1607.7368
$ws.Range("M31").Value = -595.75
$ws.Range("N31").Value = -2197.7368
$ws.Range("H34").Value = 1395.2963
$ws.Range("I34").Value = 890.75
$ws.Range("J34").Value = 1607.7368
$ws.Range("K34").Value = 890.75
$ws.Range("L34").Value = 1607.7368
$ws.Range("M34").Value = -688.75
$ws.Range("N34").Value = -2011.7368
$ws.Range("H58").Value = 2175293.5
$ws.Range("I58").Value = 3953845
$ws.Range("K58").Value = 3953845
$ws.Range("M58").Value = -3953642
$ws.Range("H68").Value = 43249.75
$ws.Range("J68").Value = 43249.75
$ws.Range("L68").Value = 43249.75
$ws.Range("N68").Value = -44747.75
$ws.Range("H71").Value = 43249.75
$ws.Range("J71").Value = 43249.75
$ws.Range("L71").Value = 129749.25
$ws.Range("N71").Value = -137237.25
$ws.Range("H134").Value = 1598.8
$ws.Range("I134").Value = 1241.0667
$ws.Range("J134").Value = 2672
$ws.Range("K134").Value = 3723.2001
$ws.Range("L134").Value = 8016
$ws.Range("M134").Value = -1188.2001
$ws.Range("N134").Value = -13086
$ws.Range("H136").Value = 2175293.5
$ws.Range("I136").Value = 3953845
$ws.Range("K136").Value = 11861535
$ws.Range("M136").Value = -11858985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4666.6665
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730
$ws.Range("H73").Value = 4666.6665
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064
$ws.Range("H80").Value = 8134.0527
$ws.Range("I80").Value = 7061
$ws.Range("J80").Value = 10459
$ws.Range("K80").Value = 7061
$ws.Range("L80").Value = 10459
$ws.Range("M80").Value = -6063
$ws.Range("N80").Value = -12455
$ws.Range("H83").Value = 8134.0527
$ws.Range("I83").Value = 7061
$ws.Range("J83").Value = 10459
$ws.Range("K83").Value = 35305
$ws.Range("L83").Value = 52295
$ws.Range("M83").Value = -30313
$ws.Range("N83").Value = -62279
$ws.Range("H97").Value = 1183.65
$ws.Range("I97").Value = 880
$ws.Range("J97").Value = 1554.7778
$ws.Range("K97").Value = 880
$ws.Range("L97").Value = 1554.7778
$ws.Range("M97").Value = -384
$ws.Range("N97").Value = -2546.7778
$ws.Range("H107").Value = 100
$ws.Range("I107").Value = 100
$ws.Range("K107").Value = 100
$ws.Range("M107").Value = 1820
$ws.Range("H122").Value = 2323.7334
$ws.Range("I122").Value = 1939.4
$ws.Range("K122").Value = 5818.200000000001
$ws.Range("M122").Value = -3368.200000000001
$ws.Range("H132").Value = 1605194.8
$ws.Range("I132").Value = 2566091
$ws.Range("K132").Value = 7698273
$ws.Range("M132").Value = -7695743

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5258.4375
$ws.Range("I16").Value = 6771.1816
$ws.Range("J16").Value = 1930.4
$ws.Range("K16").Value = 6771.1816
$ws.Range("L16").Value = 1930.4
$ws.Range("M16").Value = -6601.1816
$ws.Range("N16").Value = -2270.4
$ws.Range("H61").Value = 4311.125
$ws.Range("I61").Value = 3167
$ws.Range("J61").Value = 4997.6
$ws.Range("K61").Value = 3167
$ws.Range("L61").Value = 4997.6
$ws.Range("M61").Value = -2965
$ws.Range("N61").Value = -5401.6
$ws.Range("H113").Value = 4311.125
$ws.Range("I113").Value = 3167
$ws.Range("J113").Value = 4997.6
$ws.Range("K113").Value = 3167
$ws.Range("L113").Value = 4997.6
$ws.Range("M113").Value = -997
$ws.Range("N113").Value = -9337.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 39999.953
$ws.Range("J125").Value = 39999.953
$ws.Range("L125").Value = 39999.953
$ws.Range("N125").Value = -49839.953
$ws.Range("H126").Value = 9294.105
$ws.Range("I126").Value = 11943.909
$ws.Range("J126").Value = 5650.625
$ws.Range("K126").Value = 35831.727
$ws.Range("L126").Value = 16951.875
$ws.Range("M126").Value = -33361.727
$ws.Range("N126").Value = -21891.875
$ws.Range("H132").Value = 1458.2041
$ws.Range("I132").Value = 1076.091
$ws.Range("J132").Value = 2246.3125
$ws.Range("K132").Value = 3228.273
$ws.Range("L132").Value = 6738.9375
$ws.Range("M132").Value = -698.2729999999997
$ws.Range("N132").Value = -11798.9375
$ws.Range("H135").Value = 75423.25
$ws.Range("J135").Value = 75423.25
$ws.Range("L135").Value = 75423.25
$ws.Range("N135").Value = -85563.25
